# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" for
# the most-recently-handed-off file (c700ca98-f34d-40d1-8ec8-05c12b21372e)
# across all three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 is the c700ca98-... file; column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-08-26 14:52:53"

# zh-cn sheet: row 7 is the c700ca98-... file; column H = "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-08-26 14:52:48"

# de-de sheet: row 7 is the c700ca98-... file; column H = "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-08-26 14:52:53"
